$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.023556
$ws.Range("H2").Value = 0.07066799999999999
$ws.Range("I2").Value = 0.08088674619362546
$ws.Range("J2").Value = 0.08088674619362549
$ws.Range("M2").Value = 0.007418999999999999
$ws.Range("N2").Value = 0.022257
$ws.Range("O2").Value = 0.00543403025741878
$ws.Range("P2").Value = 0.00543403025741878
$ws.Range("Q2").Value = 0.000174761964
$ws.Range("R2").Value = 0.001572857676
$ws.Range("S2").Value = 0.0004395410262403141
$ws.Range("T2").Value = 0.0004395410262403142
$ws.Range("G3").Value = 0.023556
$ws.Range("H3").Value = 0.07066799999999999
$ws.Range("I3").Value = 0.08088674619362546
$ws.Range("J3").Value = 0.08088674619362549
$ws.Range("O3").Value = 0.3642834035872852
$ws.Range("P3").Value = 0.3642834035872852
$ws.Range("Q3").Value = 0.011715592304
$ws.Range("R3").Value = 0.105440330736
$ws.Range("S3").Value = 0.02946569920851477
$ws.Range("T3").Value = 0.02946569920851478
$ws.Range("G4").Value = 0.023556
$ws.Range("H4").Value = 0.07066799999999999
$ws.Range("I4").Value = 0.08088674619362546
$ws.Range("J4").Value = 0.08088674619362549
$ws.Range("M4").Value = 0.8605153333333333
$ws.Range("N4").Value = 2.581546
$ws.Range("O4").Value = 0.630282566155296
$ws.Range("P4").Value = 0.630282566155296
$ws.Range("Q4").Value = 0.020270299192
$ws.Range("R4").Value = 0.182432692728
$ws.Range("S4").Value = 0.05098150595887038
$ws.Range("T4").Value = 0.05098150595887039
$ws.Range("I5").Value = 0.8589438069010353
$ws.Range("J5").Value = 0.8589438069010354
$ws.Range("M5").Value = 0.007418999999999999
$ws.Range("N5").Value = 0.022257
$ws.Range("O5").Value = 0.00543403025741878
$ws.Range("P5").Value = 0.00543403025741878
$ws.Range("Q5").Value = 0.00185581339
$ws.Range("R5").Value = 0.01670232051
$ws.Range("S5").Value = 0.0046675266361227
$ws.Range("T5").Value = 0.0046675266361227
$ws.Range("I6").Value = 0.8589438069010353
$ws.Range("J6").Value = 0.8589438069010354
$ws.Range("O6").Value = 0.3642834035872852
$ws.Range("P6").Value = 0.3642834035872852
$ws.Range("S6").Value = 0.3128989734681291
$ws.Range("T6").Value = 0.3128989734681291
$ws.Range("I7").Value = 0.8589438069010353
$ws.Range("J7").Value = 0.8589438069010354
$ws.Range("M7").Value = 0.8605153333333333
$ws.Range("N7").Value = 2.581546
$ws.Range("O7").Value = 0.630282566155296
$ws.Range("P7").Value = 0.630282566155296
$ws.Range("Q7").Value = 0.2152521738644444
$ws.Range("R7").Value = 1.93726956478
$ws.Range("S7").Value = 0.5413773067967835
$ws.Range("T7").Value = 0.5413773067967836
$ws.Range("G8").Value = 0.01752266666666667
$ws.Range("H8").Value = 0.052568
$ws.Range("I8").Value = 0.0601694469053391
$ws.Range("J8").Value = 0.06016944690533912
$ws.Range("M8").Value = 0.007418999999999999
$ws.Range("N8").Value = 0.022257
$ws.Range("O8").Value = 0.00543403025741878
$ws.Range("P8").Value = 0.00543403025741878
$ws.Range("Q8").Value = 0.000130000664
$ws.Range("R8").Value = 0.001170005976
$ws.Range("S8").Value = 0.0003269625950557654
$ws.Range("T8").Value = 0.0003269625950557655
$ws.Range("G9").Value = 0.01752266666666667
$ws.Range("H9").Value = 0.052568
$ws.Range("I9").Value = 0.0601694469053391
$ws.Range("J9").Value = 0.06016944690533912
$ws.Range("O9").Value = 0.3642834035872852
$ws.Range("P9").Value = 0.3642834035872852
$ws.Range("Q9").Value = 0.008714909948444444
$ws.Range("R9").Value = 0.07843418953600001
$ws.Range("S9").Value = 0.02191873091064137
$ws.Range("T9").Value = 0.02191873091064138
$ws.Range("G10").Value = 0.01752266666666667
$ws.Range("H10").Value = 0.052568
$ws.Range("I10").Value = 0.0601694469053391
$ws.Range("J10").Value = 0.06016944690533912
$ws.Range("M10").Value = 0.8605153333333333
$ws.Range("N10").Value = 2.581546
$ws.Range("O10").Value = 0.630282566155296
$ws.Range("P10").Value = 0.630282566155296
$ws.Range("Q10").Value = 0.01507852334755555
$ws.Range("R10").Value = 0.135706710128
$ws.Range("S10").Value = 0.03792375339964196
$ws.Range("T10").Value = 0.03792375339964197
